# REPORTGEN-708: Add No Header option for components in excel reports
#
# This script adds documentation rows (for the new "HEADER=NO" option, and for
# the already-existing "EVOLUTION" option that was missing from one sheet) to
# several component-documentation sheets, and marks the new "NEW" label next
# to the component title on a couple of sheets, then changes which sheet/cell
# is active/selected.

$wb = $excel.ActiveWorkbook

$textEvolution = "* EVOLUTION=true|false to display added and removed violations columns. By default or if not exists, is true if there is a previous snapshot."
$textHeader    = "* HEADER=NO to not display headers (useful for excel report when you want to define your own customized headers). By default if option is not present or different from NO, headers are displayed"
$textNew       = "NEW"

# ---------------------------------------------------------------------------
# "3 - Evolution of standards": insert two new documentation rows (EVOLUTION
# and HEADER) right after the existing options, before the blank separator row.
# ---------------------------------------------------------------------------
$wsEvoStd = $wb.Worksheets.Item("3 - Evolution of standards")
$wsEvoStd.Rows.Item(6).Insert()
$wsEvoStd.Range("B6").Value = $textEvolution
$wsEvoStd.Rows.Item(7).Insert()
$wsEvoStd.Range("B7").Value = $textHeader
$wsEvoStd.Range("B7").Select()

# ---------------------------------------------------------------------------
# "3-CastRulesEvoByQualityCategory": insert the HEADER documentation row
# before the blank separator row.
# ---------------------------------------------------------------------------
$wsCastRules = $wb.Worksheets.Item("3-CastRulesEvoByQualityCategory")
$wsCastRules.Rows.Item(6).Insert()
$wsCastRules.Range("B6").Value = $textHeader

# ---------------------------------------------------------------------------
# "3-ListViolBookmarks": flag the component as NEW, and fill in the
# previously blank row with the HEADER documentation.
# ---------------------------------------------------------------------------
$wsListViolBk = $wb.Worksheets.Item("3-ListViolBookmarks")
$wsListViolBk.Range("C1").Value = $textNew
$wsListViolBk.Range("B7").Copy() | Out-Null
$wsListViolBk.Range("B6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$wsListViolBk.Range("B6").Value = $textHeader
$wsListViolBk.Range("B6").Select()

# ---------------------------------------------------------------------------
# "3 - List of violations statist ": insert the HEADER documentation row
# before the blank separator row, then make this sheet the active tab.
# ---------------------------------------------------------------------------
$wsListViolStat = $wb.Worksheets.Item("3 - List of violations statist ")
$wsListViolStat.Rows.Item(10).Insert()
$wsListViolStat.Range("B10").Value = $textHeader
$wsListViolStat.Activate()
$wsListViolStat.Range("B10").Select()
